$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-looking numeric strings to stay as text (Test ID, Mobile No.)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("K2").NumberFormat = "@"

# Update row 2 with the new test record values
$ws.Range("A2").Value = 58
$ws.Range("B2").Value = "26"
$ws.Range("C2").Value = "15-03-2024"
$ws.Range("D2").Value = 25.69
$ws.Range("E2").Value = 226.56
$ws.Range("F2").Value = "asdkjklasjdl"
$ws.Range("G2").Value = 25
$ws.Range("H2").Value = "Male"
$ws.Range("I2").Value = 32
$ws.Range("J2").Value = "asdkjasld"
$ws.Range("K2").Value = "9896532265"
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 20
$ws.Range("O2").Value = 300
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 49
$ws.Range("R2").Value = 25
$ws.Range("S2").Value = 30
$ws.Range("T2").Value = 0.47298236733493
$ws.Range("U2").Value = "Millets(Pearl Millet, Sorghum), Maize, Soybean, Groundnut"

# Remove the now-duplicate third record entirely
$ws.Rows.Item(3).Delete()
